$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: update title and link
$ws.Range("D9").Value = "AI가 외계 생명체를 최초로 발견하게 될까요?"
$ws.Range("E9").Value = "https://pdsi.pabii.com/ai%ea%b0%80-%ec%99%b8%ea%b3%84-%ec%83%9d%eb%aa%85%ec%b2%b4%eb%a5%bc-%ec%b5%9c%ec%b4%88%eb%a1%9c-%eb%b0%9c%ea%b2%ac%ed%95%98%ea%b2%8c-%eb%90%a0%ea%b9%8c%ec%9a%94/#utm_source=rss&utm_medium=rss&utm_campaign=ai%25ea%25b0%2580-%25ec%2599%25b8%25ea%25b3%2584-%25ec%2583%259d%25eb%25aa%2585%25ec%25b2%25b4%25eb%25a5%25bc-%25ec%25b5%259c%25ec%25b4%2588%25eb%25a1%259c-%25eb%25b0%259c%25ea%25b2%25ac%25ed%2595%2598%25ea%25b2%258c-%25eb%2590%25a0%25ea%25b9%258c%25ec%259a%2594"

# Row 36: update title and link
$ws.Range("D36").Value = "Semi-supervised Learning of FixMatch and after FixMatch"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/395"

# Row 50: update title and link
$ws.Range("D50").Value = "correlation"
$ws.Range("E50").Value = "http://incredible.egloos.com/7578543"
